# Update "想去人数" (F4/F5/F6) figures on both the "展览" sheet and the
# "全部类型" sheet to reflect the latest generated numbers.
$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 31
    $ws.Range("F5").Value = 2345
    $ws.Range("F6").Value = 217
}
